{"js": "// Cinematic camera, end prototype Day 1\n// Applies the text edits described by the diff to the Word document body.\nconst body = context.document.body;\n\n// ---------------------------------------------------------------\n// 1) Typo fix in the intro paragraph: remove the double space before\n//    \"\u0430 \u043d\u043e\u0447\u044c\u044e\" and fix \"\u043f\u0440\u044f\u0447\u044e\u0442\u0441\u044f\" -> \"\u043f\u0440\u044f\u0447\u0443\u0442\u0441\u044f\".\n// ---------------------------------------------------------------\nconst introHit = body.search(\n  \"\u0432\u044b\u0445\u043e\u0434\u044f\u0442,  \u0430 \u043d\u043e\u0447\u044c\u044e \u0436\u0438\u0442\u0435\u043b\u0438 \u043f\u0440\u044f\u0447\u044e\u0442\u0441\u044f \u043f\u043e \u0434\u043e\u043c\u0430\u043c\",\n  { matchCase: true }\n);\nintroHit.load(\"text\");\nawait context.sync();\nif (introHit.items.length > 0) {\n  introHit.items[0].insertText(\n    \"\u0432\u044b\u0445\u043e\u0434\u044f\u0442, \u0430 \u043d\u043e\u0447\u044c\u044e \u0436\u0438\u0442\u0435\u043b\u0438 \u043f\u0440\u044f\u0447\u0443\u0442\u0441\u044f \u043f\u043e \u0434\u043e\u043c\u0430\u043c\",\n    \"Replace\"\n  );\n}\n\n// ---------------------------------------------------------------\n// 2) \"\u0418\u0433\u0440\u0430 \u043d\u0430\u0447\u0438\u043d\u0430\u0435\u0442\u0441\u044f...\" paragraph (Day 1): insert the new ghost\n//    scene between \"\u043f\u043e\u0442\u043e\u043c \u0432 \u043a\u0430\u043a\u043e\u0439-\u0442\u043e \u043c\u043e\u043c\u0435\u043d\u0442\" and the old\n//    \"\u0432\u044b\u0431\u0435\u0433\u0430\u0435\u0442 \u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u044e\" text.\n// ---------------------------------------------------------------\nconst day1Hit = body.search(\n  \"\u043f\u043e\u0442\u043e\u043c \u0432 \u043a\u0430\u043a\u043e\u0439-\u0442\u043e \u043c\u043e\u043c\u0435\u043d\u0442 \u0432\u044b\u0431\u0435\u0433\u0430\u0435\u0442 \u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u044e\",\n  { matchCase: true }\n);\nday1Hit.load(\"text\");\nawait context.sync();\nif (day1Hit.items.length > 0) {\n  day1Hit.items[0].insertText(\n    \"\u043f\u043e\u0442\u043e\u043c \u0432 \u043a\u0430\u043a\u043e\u0439-\u0442\u043e \u043c\u043e\u043c\u0435\u043d\u0442 \u0432\u0438\u0434\u0438\u043c \u043f\u0440\u0438\u0437\u0440\u0430\u043a\u0430 \u0438 \u0431\u0435\u0436\u0438\u043c \u0437\u0430 \u043d\u0438\u043c, \u0433\u043b\u0430\u0432\u043d\u044b\u0439 \u0433\u0435\u0440\u043e\u0439 \u0433\u043e\u0432\u043e\u0440\u0438\u0442 \u043c\u043d\u0435 \u043d\u0443\u0436\u043d\u044b \u043e\u0442\u0432\u0435\u0442\u0430, \u043f\u0440\u0438\u0437\u0440\u0430\u043a \u0433\u043e\u0432\u043e\u0440\u0438\u0442 \u0447\u0442\u043e \u0442\u0435\u0431\u0435 \u0441\u0442\u043e\u0438\u0442 \u043f\u043e\u0433\u043e\u0432\u043e\u0440\u0438\u0442\u044c \u0441 \u043a\u0443\u0437\u043d\u0435\u0446\u043e\u043c, \u043e\u043d \u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u0435 \u043d\u0430\u043f\u0440\u043e\u0442\u0438\u0432, \u0433\u043b\u0430\u0432\u043d\u044b\u0439 \u0433\u0435\u0440\u043e\u0439 \u0432\u0431\u0435\u0433\u0430\u0435\u0442 \u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u044e\",\n    \"Replace\"\n  );\n}\n\n// ---------------------------------------------------------------\n// 3) \"\u0414\u0435\u043d\u044c 2\" paragraph:\n//    a) insert \"\u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u0435 \u0438\u0434\u0435\u0442 \u0441\u043d\u0435\u0433,\" after \"\u0436\u0438\u0442\u0435\u043b\u0438 \u0435\u0433\u043e \u0431\u043b\u0430\u0433\u043e\u0434\u0430\u0440\u044f\u0442,\"\n//    b) drop \", \u0447\u0442\u043e \u0442\u0443\u0442 \u043f\u0440\u043e\u0438\u0441\u0445\u043e\u0434\u0438\u0442, \u0438 \u0433\u0434\u0435 \u043e\u043d, \u043e\u043d\u0438 \u0435\u043c\u0443 \u043d\u0435 \u043e\u0442\u0432\u0435\u0447\u0430\u044e\u0442\"\n//    c) append the closing voice-over line at the very end of the\n//       paragraph.\n// ---------------------------------------------------------------\nconst day2HitA = body.search(\n  \"\u0436\u0438\u0442\u0435\u043b\u0438 \u0435\u0433\u043e \u0431\u043b\u0430\u0433\u043e\u0434\u0430\u0440\u044f\u0442, \u043e\u043d \u0445\u043e\u0434\u0438\u0442 \u043f\u043e \u0434\u0435\u0440\u0435\u0432\u043d\u0438, \u0438 \u0440\u0430\u0441\u0441\u043f\u0440\u0430\u0448\u0438\u0432\u0430\u0435\u0442 \u0436\u0438\u0442\u0435\u043b\u0435 \u0433\u0434\u0435 \u043a\u0443\u0437\u043d\u0435\u0446, \u0447\u0442\u043e \u0442\u0443\u0442 \u043f\u0440\u043e\u0438\u0441\u0445\u043e\u0434\u0438\u0442, \u0438 \u0433\u0434\u0435 \u043e\u043d, \u043e\u043d\u0438 \u0435\u043c\u0443 \u043d\u0435 \u043e\u0442\u0432\u0435\u0447\u0430\u044e\u0442,\",\n  { matchCase: true }\n);\nday2HitA.load(\"text\");\nawait context.sync();\nif (day2HitA.items.length > 0) {\n  day2HitA.items[0].insertText(\n    \"\u0436\u0438\u0442\u0435\u043b\u0438 \u0435\u0433\u043e \u0431\u043b\u0430\u0433\u043e\u0434\u0430\u0440\u044f\u0442, \u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u0435 \u0438\u0434\u0435\u0442 \u0441\u043d\u0435\u0433, \u043e\u043d \u0445\u043e\u0434\u0438\u0442 \u043f\u043e \u0434\u0435\u0440\u0435\u0432\u043d\u0438, \u0438 \u0440\u0430\u0441\u0441\u043f\u0440\u0430\u0448\u0438\u0432\u0430\u0435\u0442 \u0436\u0438\u0442\u0435\u043b\u0435 \u0433\u0434\u0435 \u043a\u0443\u0437\u043d\u0435\u0446, \u0433\u043b\u0430\u0432\u043d\u044b\u0439 \u0433\u0435\u0440\u043e\u0439 \u043d\u0430\u0445\u043e\u0434\u0438\u0442 \u0437\u0430\u043f\u0438\u0441\u043a\u0443 \u043e\u0442 \u043a\u0443\u0437\u043d\u0435\u0446\u0430, \u0447\u0442\u043e \u043e\u043d \u043f\u043e\u0448\u0435\u043b \u043d\u0430 \u0440\u044b\u0431\u0430\u043b\u043a\u0443,\",\n    \"Replace\"\n  );\n}\n\n// The original text already contained \" \u0433\u043b\u0430\u0432\u043d\u044b\u0439 \u0433\u0435\u0440\u043e\u0439 \u043d\u0430\u0445\u043e\u0434\u0438\u0442 \u0437\u0430\u043f\u0438\u0441\u043a\u0443 \u043e\u0442\n// \u043a\u0443\u0437\u043d\u0435\u0446\u0430, \u0447\u0442\u043e \u043e\u043d \u043f\u043e\u0448\u0435\u043b \u043d\u0430 \u0440\u044b\u0431\u0430\u043b\u043a\u0443,\" right after the part we just replaced;\n// remove that now-duplicated fragment.\nconst day2HitB = body.search(\n  \"\u0440\u044b\u0431\u0430\u043b\u043a\u0443, \u0433\u043b\u0430\u0432\u043d\u044b\u0439 \u0433\u0435\u0440\u043e\u0439 \u043d\u0430\u0445\u043e\u0434\u0438\u0442 \u0437\u0430\u043f\u0438\u0441\u043a\u0443 \u043e\u0442 \u043a\u0443\u0437\u043d\u0435\u0446\u0430, \u0447\u0442\u043e \u043e\u043d \u043f\u043e\u0448\u0435\u043b \u043d\u0430 \u0440\u044b\u0431\u0430\u043b\u043a\u0443,\",\n  { matchCase: true }\n);\nday2HitB.load(\"text\");\nawait context.sync();\nif (day2HitB.items.length > 0) {\n  day2HitB.items[0].insertText(\"\u0440\u044b\u0431\u0430\u043b\u043a\u0443,\", \"Replace\");\n}\n\nconst day2HitC = body.search(\n  \"\u043f\u043e\u0442\u043e\u043c\u0443 \u0447\u0442\u043e \u0442\u0430\u043c \u043f\u0440\u043e\u0438\u0441\u0445\u043e\u0434\u0438\u0442\u044c \u0447\u0442\u043e-\u0442\u043e \u0441\u0442\u0440\u0430\u043d\u043d\u043e\u0435\u2026\",\n  { matchCase: true }\n);\nday2HitC.load(\"text\");\nawait context.sync();\nif (day2HitC.items.length > 0) {\n  day2HitC.items[0].insertText(\n    \"\u043f\u043e\u0442\u043e\u043c\u0443 \u0447\u0442\u043e \u0442\u0430\u043c \u043f\u0440\u043e\u0438\u0441\u0445\u043e\u0434\u0438\u0442\u044c \u0447\u0442\u043e-\u0442\u043e \u0441\u0442\u0440\u0430\u043d\u043d\u043e\u0435, \u0413\u043e\u043b\u043e\u0441 \u0437\u0430 \u043a\u0430\u0434\u0440\u043e\u043c: \u201e\u0412\u043e\u0437\u043c\u043e\u0436\u043d\u043e \u044d\u0442\u043e \u043d\u0430\u0447\u0430\u043b\u043e \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u043d\u043e\u0439 \u0434\u0440\u0443\u0436\u0431\u044b, \u043a\u0442\u043e \u0437\u043d\u0430\u0435\u0442\u2026\u201d\",\n    \"Replace\"\n  );\n}\n\n// ---------------------------------------------------------------\n// 4) Remove the extra blank paragraphs:\n//    - two of the three blank paragraphs right after\n//      \"\u0422\u0435\u0441\u0442\u043e\u0432\u043e\u0435 \u043d\u0430\u0437\u0432\u0430\u043d\u0438\u0435 \u0438\u0433\u0440\u044b...\"\n//    - the trailing blank paragraph at the very end of the document\n// ---------------------------------------------------------------\nlet paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nlet items = paras.items;\n\n// 4a) Two of the three blank paragraphs following \"\u0422\u0435\u0441\u0442\u043e\u0432\u043e\u0435 \u043d\u0430\u0437\u0432\u0430\u043d\u0438\u0435\n// \u0438\u0433\u0440\u044b...\": delete the first two of that run (keep the third, closest\n// to the next piece of text).\nconst blankRunAfterTitle = [];\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"\u0422\u0435\u0441\u0442\u043e\u0432\u043e\u0435 \u043d\u0430\u0437\u0432\u0430\u043d\u0438\u0435 \u0438\u0433\u0440\u044b\") !== -1) {\n    let j = i + 1;\n    while (j < items.length && items[j].text === \"\") {\n      blankRunAfterTitle.push(j);\n      j++;\n    }\n    break;\n  }\n}\nif (blankRunAfterTitle.length >= 2) {\n  // Delete from the highest index down so earlier indices stay valid.\n  items[blankRunAfterTitle[1]].delete();\n  items[blankRunAfterTitle[0]].delete();\n  await context.sync();\n}\n\n// 4b) The trailing blank paragraph at the very end of the document is\n// the body's terminal paragraph, so a plain Paragraph.delete() is a\n// no-op for it. Instead, expand a range from the end of the\n// second-to-last paragraph through the end of the body and delete\n// that range, which removes the paragraph mark as well.\nparas = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\nitems = paras.items;\nif (items.length >= 2 && items[items.length - 1].text === \"\") {\n  const secondLast = items[items.length - 2];\n  const secondLastEnd = secondLast.getRange(\"End\");\n  const bodyEnd = body.getRange(\"End\");\n  const trailingRange = secondLastEnd.expandTo(bodyEnd);\n  trailingRange.delete();\n  await context.sync();\n}\n", "ps1": "# Cinematic camera, end prototype Day 1\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------\n# 1) Typo fix in the intro paragraph: remove the double space before\n#    \"\u0430 \u043d\u043e\u0447\u044c\u044e\" and fix \"\u043f\u0440\u044f\u0447\u044e\u0442\u0441\u044f\" -> \"\u043f\u0440\u044f\u0447\u0443\u0442\u0441\u044f\".\n# ---------------------------------------------------------------\n$d.Content.Find.Execute(\n    \"\u0432\u044b\u0445\u043e\u0434\u044f\u0442,  \u0430 \u043d\u043e\u0447\u044c\u044e \u0436\u0438\u0442\u0435\u043b\u0438 \u043f\u0440\u044f\u0447\u044e\u0442\u0441\u044f \u043f\u043e \u0434\u043e\u043c\u0430\u043c\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"\u0432\u044b\u0445\u043e\u0434\u044f\u0442, \u0430 \u043d\u043e\u0447\u044c\u044e \u0436\u0438\u0442\u0435\u043b\u0438 \u043f\u0440\u044f\u0447\u0443\u0442\u0441\u044f \u043f\u043e \u0434\u043e\u043c\u0430\u043c\", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 2) \"\u0418\u0433\u0440\u0430 \u043d\u0430\u0447\u0438\u043d\u0430\u0435\u0442\u0441\u044f...\" paragraph (Day 1): insert the new ghost\n#    scene between \"\u043f\u043e\u0442\u043e\u043c \u0432 \u043a\u0430\u043a\u043e\u0439-\u0442\u043e \u043c\u043e\u043c\u0435\u043d\u0442\" and the old\n#    \"\u0432\u044b\u0431\u0435\u0433\u0430\u0435\u0442 \u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u044e\" text.\n# ---------------------------------------------------------------\n$d.Content.Find.Execute(\n    \"\u043f\u043e\u0442\u043e\u043c \u0432 \u043a\u0430\u043a\u043e\u0439-\u0442\u043e \u043c\u043e\u043c\u0435\u043d\u0442 \u0432\u044b\u0431\u0435\u0433\u0430\u0435\u0442 \u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u044e\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"\u043f\u043e\u0442\u043e\u043c \u0432 \u043a\u0430\u043a\u043e\u0439-\u0442\u043e \u043c\u043e\u043c\u0435\u043d\u0442 \u0432\u0438\u0434\u0438\u043c \u043f\u0440\u0438\u0437\u0440\u0430\u043a\u0430 \u0438 \u0431\u0435\u0436\u0438\u043c \u0437\u0430 \u043d\u0438\u043c, \u0433\u043b\u0430\u0432\u043d\u044b\u0439 \u0433\u0435\u0440\u043e\u0439 \u0433\u043e\u0432\u043e\u0440\u0438\u0442 \u043c\u043d\u0435 \u043d\u0443\u0436\u043d\u044b \u043e\u0442\u0432\u0435\u0442\u0430, \u043f\u0440\u0438\u0437\u0440\u0430\u043a \u0433\u043e\u0432\u043e\u0440\u0438\u0442 \u0447\u0442\u043e \u0442\u0435\u0431\u0435 \u0441\u0442\u043e\u0438\u0442 \u043f\u043e\u0433\u043e\u0432\u043e\u0440\u0438\u0442\u044c \u0441 \u043a\u0443\u0437\u043d\u0435\u0446\u043e\u043c, \u043e\u043d \u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u0435 \u043d\u0430\u043f\u0440\u043e\u0442\u0438\u0432, \u0433\u043b\u0430\u0432\u043d\u044b\u0439 \u0433\u0435\u0440\u043e\u0439 \u0432\u0431\u0435\u0433\u0430\u0435\u0442 \u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u044e\", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 3) \"\u0414\u0435\u043d\u044c 2\" paragraph:\n#    a) insert \"\u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u0435 \u0438\u0434\u0435\u0442 \u0441\u043d\u0435\u0433,\" after \"\u0436\u0438\u0442\u0435\u043b\u0438 \u0435\u0433\u043e \u0431\u043b\u0430\u0433\u043e\u0434\u0430\u0440\u044f\u0442,\"\n#       and drop \", \u0447\u0442\u043e \u0442\u0443\u0442 \u043f\u0440\u043e\u0438\u0441\u0445\u043e\u0434\u0438\u0442, \u0438 \u0433\u0434\u0435 \u043e\u043d, \u043e\u043d\u0438 \u0435\u043c\u0443 \u043d\u0435 \u043e\u0442\u0432\u0435\u0447\u0430\u044e\u0442\"\n#    b) remove the now-duplicated \"\u0433\u043b\u0430\u0432\u043d\u044b\u0439 \u0433\u0435\u0440\u043e\u0439 \u043d\u0430\u0445\u043e\u0434\u0438\u0442 \u0437\u0430\u043f\u0438\u0441\u043a\u0443...\"\n#       fragment\n#    c) append the closing voice-over line at the very end.\n# ---------------------------------------------------------------\n$d.Content.Find.Execute(\n    \"\u0436\u0438\u0442\u0435\u043b\u0438 \u0435\u0433\u043e \u0431\u043b\u0430\u0433\u043e\u0434\u0430\u0440\u044f\u0442, \u043e\u043d \u0445\u043e\u0434\u0438\u0442 \u043f\u043e \u0434\u0435\u0440\u0435\u0432\u043d\u0438, \u0438 \u0440\u0430\u0441\u0441\u043f\u0440\u0430\u0448\u0438\u0432\u0430\u0435\u0442 \u0436\u0438\u0442\u0435\u043b\u0435 \u0433\u0434\u0435 \u043a\u0443\u0437\u043d\u0435\u0446, \u0447\u0442\u043e \u0442\u0443\u0442 \u043f\u0440\u043e\u0438\u0441\u0445\u043e\u0434\u0438\u0442, \u0438 \u0433\u0434\u0435 \u043e\u043d, \u043e\u043d\u0438 \u0435\u043c\u0443 \u043d\u0435 \u043e\u0442\u0432\u0435\u0447\u0430\u044e\u0442,\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"\u0436\u0438\u0442\u0435\u043b\u0438 \u0435\u0433\u043e \u0431\u043b\u0430\u0433\u043e\u0434\u0430\u0440\u044f\u0442, \u0432 \u0434\u0435\u0440\u0435\u0432\u043d\u0435 \u0438\u0434\u0435\u0442 \u0441\u043d\u0435\u0433, \u043e\u043d \u0445\u043e\u0434\u0438\u0442 \u043f\u043e \u0434\u0435\u0440\u0435\u0432\u043d\u0438, \u0438 \u0440\u0430\u0441\u0441\u043f\u0440\u0430\u0448\u0438\u0432\u0430\u0435\u0442 \u0436\u0438\u0442\u0435\u043b\u0435 \u0433\u0434\u0435 \u043a\u0443\u0437\u043d\u0435\u0446, \u0433\u043b\u0430\u0432\u043d\u044b\u0439 \u0433\u0435\u0440\u043e\u0439 \u043d\u0430\u0445\u043e\u0434\u0438\u0442 \u0437\u0430\u043f\u0438\u0441\u043a\u0443 \u043e\u0442 \u043a\u0443\u0437\u043d\u0435\u0446\u0430, \u0447\u0442\u043e \u043e\u043d \u043f\u043e\u0448\u0435\u043b \u043d\u0430 \u0440\u044b\u0431\u0430\u043b\u043a\u0443,\", 2) | Out-Null\n\n$d.Content.Find.Execute(\n    \"\u0440\u044b\u0431\u0430\u043b\u043a\u0443, \u0433\u043b\u0430\u0432\u043d\u044b\u0439 \u0433\u0435\u0440\u043e\u0439 \u043d\u0430\u0445\u043e\u0434\u0438\u0442 \u0437\u0430\u043f\u0438\u0441\u043a\u0443 \u043e\u0442 \u043a\u0443\u0437\u043d\u0435\u0446\u0430, \u0447\u0442\u043e \u043e\u043d \u043f\u043e\u0448\u0435\u043b \u043d\u0430 \u0440\u044b\u0431\u0430\u043b\u043a\u0443,\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"\u0440\u044b\u0431\u0430\u043b\u043a\u0443,\", 2) | Out-Null\n\n$d.Content.Find.Execute(\n    \"\u043f\u043e\u0442\u043e\u043c\u0443 \u0447\u0442\u043e \u0442\u0430\u043c \u043f\u0440\u043e\u0438\u0441\u0445\u043e\u0434\u0438\u0442\u044c \u0447\u0442\u043e-\u0442\u043e \u0441\u0442\u0440\u0430\u043d\u043d\u043e\u0435\u2026\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"\u043f\u043e\u0442\u043e\u043c\u0443 \u0447\u0442\u043e \u0442\u0430\u043c \u043f\u0440\u043e\u0438\u0441\u0445\u043e\u0434\u0438\u0442\u044c \u0447\u0442\u043e-\u0442\u043e \u0441\u0442\u0440\u0430\u043d\u043d\u043e\u0435, \u0413\u043e\u043b\u043e\u0441 \u0437\u0430 \u043a\u0430\u0434\u0440\u043e\u043c: \u201e\u0412\u043e\u0437\u043c\u043e\u0436\u043d\u043e \u044d\u0442\u043e \u043d\u0430\u0447\u0430\u043b\u043e \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u043d\u043e\u0439 \u0434\u0440\u0443\u0436\u0431\u044b, \u043a\u0442\u043e \u0437\u043d\u0430\u0435\u0442\u2026\u201d\", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 4) Remove the extra blank paragraphs:\n#    - two of the three blank paragraphs right after\n#      \"\u0422\u0435\u0441\u0442\u043e\u0432\u043e\u0435 \u043d\u0430\u0437\u0432\u0430\u043d\u0438\u0435 \u0438\u0433\u0440\u044b...\"\n#    - the trailing blank paragraph at the very end of the document\n# ---------------------------------------------------------------\n$titleIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -match \"\u0422\u0435\u0441\u0442\u043e\u0432\u043e\u0435 \u043d\u0430\u0437\u0432\u0430\u043d\u0438\u0435 \u0438\u0433\u0440\u044b\") {\n        $titleIndex = $i\n        break\n    }\n}\nif ($titleIndex -gt 0) {\n    $blankIdx = @()\n    $j = $titleIndex + 1\n    while (($j -le $d.Paragraphs.Count) -and ($d.Paragraphs.Item($j).Range.Text.Trim().Length -eq 0)) {\n        $blankIdx += $j\n        $j++\n    }\n    if ($blankIdx.Count -ge 2) {\n        # Delete the first two blanks of that run; delete the higher index\n        # first so the lower index stays valid.\n        $d.Paragraphs.Item($blankIdx[1]).Range.Delete()\n        $d.Paragraphs.Item($blankIdx[0]).Range.Delete()\n    }\n}\n\n# The trailing blank paragraph is the document's terminal paragraph, so a\n# plain Range.Delete() on it is a no-op. Instead delete the paragraph\n# mark of the paragraph right before it, which merges the (empty) last\n# paragraph away and leaves the preceding content paragraph as the new\n# final paragraph.\n$lastIndex = $d.Paragraphs.Count\nif (($lastIndex -ge 2) -and ($d.Paragraphs.Item($lastIndex).Range.Text.Trim().Length -eq 0)) {\n    $prevPara = $d.Paragraphs.Item($lastIndex - 1)\n    $markStart = $prevPara.Range.End - 1\n    $markEnd = $prevPara.Range.End\n    $d.Range($markStart, $markEnd).Delete()\n}\n"}
